$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.484439164362272
$ws.Range("B2").Value = -1.27922042006967

$ws.Range("A3").Value = -0.4646261626649179
$ws.Range("B3").Value = -0.5894678462724917

$ws.Range("A4").Value = -1.002315334143556
$ws.Range("B4").Value = -0.8420532791301936

$ws.Range("A5").Value = -0.725020227433065
$ws.Range("B5").Value = -0.6808037628073165

$ws.Range("A6").Value = 0.8211783586799222
$ws.Range("B6").Value = 0.6820964583389857

$ws.Range("A7").Value = -0.09213892798072537
$ws.Range("B7").Value = 0.005658836007705776

$ws.Range("A8").Value = 0.7902620007208699
$ws.Range("B8").Value = 0.6113274057277852

$ws.Range("A9").Value = 0.3284874063919779
$ws.Range("B9").Value = 0.3285947465387778

$ws.Range("A10").Value = -0.1876080986239586
$ws.Range("B10").Value = -0.04585113813152625

$ws.Range("A11").Value = -0.2903066199609222
$ws.Range("B11").Value = -0.1285775647085816
